# Add data for 2022-09-22
# Update the "through" date in the sheet tab name and the header cell,
# then add the new carjacking counts that were folded into the historical
# monthly totals (neighborhood x month grid).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet tab name: 2022-09-13 -> 2022-09-14
$ws.Name = "Through 2022-09-14"

# Column header text for the "current" (in-progress) month.
$ws.Range("B1").Value = "September 2022 (through September 14)"

# Updated neighborhood/month counts (carjackings).
$ws.Range("B2").Value = 2
$ws.Range("T2").Value = 3
$ws.Range("T3").Value = 2
$ws.Range("BD8").Value = 2
$ws.Range("BM10").Value = 1
$ws.Range("T12").Value = 3
$ws.Range("K15").Value = 4
$ws.Range("B17").Value = 1
$ws.Range("T24").Value = 2
$ws.Range("K38").Value = 3
$ws.Range("K42").Value = 2
$ws.Range("AU50").Value = 3
$ws.Range("AL52").Value = 1
$ws.Range("BM55").Value = 1
$ws.Range("T57").Value = 2
$ws.Range("AU64").Value = 1
$ws.Range("B74").Value = 2
$ws.Range("BM76").Value = 1
$ws.Range("K78").Value = 1
$ws.Range("T95").Value = 3
$ws.Range("AU96").Value = 1
